$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) values, and fix the WrappedBTC/Chainlink row swap (B,C,D,E for rows 16-17)
$ws.Range("D2").Value = '69.979.33'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '3.551.39'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.76'
$ws.Range("E5").Value = '  -1.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '197.35'
$ws.Range("E6").Value = '  +6.32%  '
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.210'
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.659'
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.22'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("E12").Value = '  -1.83%  '
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '4.108.20'
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '600.29'
$ws.Range("E15").Value = '  -5.08%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '70.179.94'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.19'
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("D19").Value = '3.554.20'
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.25'
$ws.Range("E22").Value = '  +3.92%  '
$ws.Range("E23").Value = '  +7.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.38'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.12'
$ws.Range("E26").Value = '  +3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.98'
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.70'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.52'
$ws.Range("E30").Value = '  +22.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.76'
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.51'
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("E35").Value = '  +6.94%  '
$ws.Range("D36").Value = '3.739.79'
$ws.Range("E36").Value = '  +6.73%  '
$ws.Range("E37").Value = '  -4.35%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.64'
$ws.Range("E39").Value = '  +3.11%  '
$ws.Range("E40").Value = '  -1.46%  '
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '499.38'
$ws.Range("E42").Value = '  -5.67%  '
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("E45").Value = '  -3.67%  '
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.01'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("E49").Value = '  -5.23%  '
$ws.Range("E50").Value = '  +2.09%  '
$ws.Range("E51").Value = '  +11.76%  '
# Some new D-column values are plain decimals (e.g. "603.76"); Excel auto-converts such
# strings to numbers when assigned directly, which changes both the stored type and
# introduces floating point artifacts. We forced those cells to text via NumberFormat "@"
# above; now restore their style to the workbook default ("Normal") so no stray number
# formatting is left behind on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D48").Style = "Normal"